$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The data rows 4-10 hold 7 (kind / catalog / description) records. The last
# record ("high-speed rail" - row 10) needs to become the first record (row 4),
# and every other record shifts down by one row.
#
# Formatting: row 10 used the "highlighted" style (s=5) while rows 4-9 used the
# plain style (s=4). After the move, row 4 should carry the highlighted style
# and rows 5-10 the plain style - i.e. the *styles* move together with the
# data, they don't stay pinned to a physical row.
# ---------------------------------------------------------------------------

# 1) Move the formatting first, while the donor rows still hold their
#    original look:
#    - row 4 adopts row 10's style (the "highlighted" one)
#    - row 10 adopts row 9's style (the "plain" one, still untouched here)
$ws.Range("A10:C10").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# 2) Now write the actual cell text, in the new order, using 2-D arrays so
#    each 3-cell row is written in a single Value assignment.
function Set-Row3($rowNum, $a, $b, $c) {
    $arr = New-Object 'object[,]' 1,3
    $arr[0,0] = $a
    $arr[0,1] = $b
    $arr[0,2] = $c
    $ws.Range("A$rowNum`:C$rowNum").Value = $arr
}

Set-Row3 4  'name%=%"高速"||name%=%"城际"||name%=%"客运专线"' 'catalog=05010201' '高铁'
Set-Row3 5  'kind=0242' 'catalog=050102' '铁路无属性'
Set-Row3 6  'kind=0243' 'catalog=050103' '铁路隧道'
Set-Row3 7  'kind=0244' 'catalog=050302' '磁悬浮无属性'
Set-Row3 8  'kind=0241' 'catalog=050303' '磁悬浮隧道'
Set-Row3 9  'kind=0245' 'catalog=050202' '地铁、轻轨无属性'
Set-Row3 10 'kind=0246' 'catalog=050203' '地铁、轻轨隧道'

# 3) The saved selection moves from A11:XFD97 to a single cell, B12.
[void]$ws.Range("B12").Select()
